# FALECPV-RecEmitidas.xlsx edit
#
# 1) Remove the "AUT S.R.I." column entirely (header + data column F) -
#    deleting the whole column shifts every later column one slot to the
#    left (G->F, H->G, ... M->L) and drops the now-unused "AUT S.R.I."
#    shared string, matching the sheet1.xml / sharedStrings.xml diff
#    (dimension A2:M10 -> A2:L10, mergeCell A2:M2 -> A2:L2, col widths
#    shift left by one, etc.)
# 2) Rename the report title from "RECIBOS EMITIDAS" to "RECIBOS EMITIDOS".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire column F ("AUT S.R.I.") - everything to its right shifts left.
$ws.Columns("F:F").Delete()

# Fix title wording (cell A2, merged across the row).
$ws.Range("A2").Value = "RECIBOS EMITIDOS"
